$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B to make room for the two new
# weekly report columns (Jun_17 and Jun_15). This shifts the existing
# "Jun_13" column (B) to D and the existing "Jun_10" column (C) to E.
$ws.Columns("B:C").Insert()

# New header row values (most-recent week first)
$ws.Range("B1").Value2 = "Jun_17"
$ws.Range("C1").Value2 = "Jun_15"

# Match the column widths used by the rest of the report (~8 chars raw
# OOXML width, i.e. 7.1667 in Excel's ColumnWidth units) for the three
# data columns C, D, E.
$ws.Columns("C").ColumnWidth = 7.1666666666
$ws.Columns("D").ColumnWidth = 7.1666666666
$ws.Columns("E").ColumnWidth = 7.1666666666

# Default every analyst row to "UN" (unchanged) for both new week
# columns; most firms did not issue a new rating this period.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value2 = "UN"
    $ws.Cells.Item($r, 3).Value2 = "UN"
}

# Wells Fargo & Co (row 21) raised its target on 6/14/2018 during the
# Jun_15 reporting week - record the rating change and highlight the
# cell the same way the other rating-change cells are highlighted.
$ws.Range("C21").Value2 = "6/14/2018,Raises Target,Market Perform -> Underperform,`$42.00 -> `$60.00"
$ws.Range("C21").Interior.Color = 13434828

Write-Output "Applied Jun_15/Jun_17 columns and Wells Fargo & Co update"
